$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all Start/End date text values in columns B and C (rows 2-16) forward by one year.
$ws.Range("B2").Value = "2024-12-01"
$ws.Range("C2").Value = "2025-01-31"
$ws.Range("B3").Value = "2025-03-15"
$ws.Range("C3").Value = "2025-03-31"
$ws.Range("B4").Value = "2025-01-15"
$ws.Range("C4").Value = "2025-01-31"
$ws.Range("B5").Value = "2025-04-01"
$ws.Range("C5").Value = "2025-04-15"
$ws.Range("B6").Value = "2025-01-01"
$ws.Range("C6").Value = "2025-01-14"
$ws.Range("B7").Value = "2025-01-15"
$ws.Range("C7").Value = "2025-01-31"
$ws.Range("B8").Value = "2025-03-16"
$ws.Range("C8").Value = "2025-03-31"
$ws.Range("B9").Value = "2025-03-24"
$ws.Range("C9").Value = "2025-04-21"
$ws.Range("B10").Value = "2025-04-25"
$ws.Range("C10").Value = "2025-06-15"
$ws.Range("B11").Value = "2024-12-01"
$ws.Range("C11").Value = "2024-12-23"
$ws.Range("B12").Value = "2025-01-01"
$ws.Range("C12").Value = "2025-01-15"
$ws.Range("B13").Value = "2025-01-16"
$ws.Range("C13").Value = "2025-01-31"
$ws.Range("B14").Value = "2025-04-01"
$ws.Range("C14").Value = "2025-04-30"
$ws.Range("B15").Value = "2025-05-01"
$ws.Range("C15").Value = "2025-05-31"
$ws.Range("B16").Value = "2025-06-01"
$ws.Range("C16").Value = "2025-06-30"

# Update the active cell selection to match the new view (F8).
$ws.Range("F8").Select()
